# building.xlsx — "Established the license upgrades": rename shop_item -> building,
# replace type/stand/unlocked columns with a merged placeable_type column and a
# new store_area column, and restyle the data table (alignment + column widths).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet (the "building" defined name follows automatically) ---
$ws.Name = "building"

# --- id / name / value / sprite_path columns (unchanged header + row text, just
#     re-pointed at the same literal strings) ---
$ws.Cells.Item(1,1).Value = "id"
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,4).Value = "value"
$ws.Cells.Item(1,5).Value = "sprite_path"

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "Shelf"
$ws.Cells.Item(2,4).Value = 300
$ws.Cells.Item(2,5).Value = "res://Asset/Building/Shelf/simple_shelf.png"

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "Stand"
$ws.Cells.Item(3,4).Value = 250
$ws.Cells.Item(3,5).Value = "res://Asset/Building/Stand/simple_stand.png"

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "Hanger"
$ws.Cells.Item(4,4).Value = 200
$ws.Cells.Item(4,5).Value = "res://Asset/Building/Hanger/simple_hanger.png"

# --- new store_area column (F) ---
$ws.Cells.Item(1,6).Value = "store_area"
$ws.Cells.Item(2,6).Value = "Beverages"
$ws.Cells.Item(3,6).Value = "Weapons"
$ws.Cells.Item(4,6).Value = "Dairy"

# --- new placeable_type column (C), replacing type/stand ---
$ws.Cells.Item(1,3).Value = "placeable_type"
$ws.Cells.Item(3,3).Value = "Weapon"
$ws.Cells.Item(2,3).Value = "Food,Drink,Material,Craft"
$ws.Cells.Item(4,3).Value = "Meat,Medicine"

# --- Alignment: value header centered, value column right-aligned, the new
#     placeable_type column (incl. the stray formatted cell at C9) left-aligned ---
$ws.Cells.Item(1,4).HorizontalAlignment = -4108  # xlCenter
$ws.Range("D2:D4").HorizontalAlignment = -4152   # xlRight
$ws.Range("C2:C4").HorizontalAlignment = -4131   # xlLeft
$ws.Cells.Item(9,3).HorizontalAlignment = -4131  # xlLeft

# --- Column widths ---
$ws.Columns.Item(3).ColumnWidth = 26.666666666666668
$ws.Columns.Item(4).ColumnWidth = 8.166666666666666
$ws.Columns.Item(5).ColumnWidth = 43.666666666666664
$ws.Columns.Item(6).ColumnWidth = 10.5

# --- Final selection ---
$ws.Range("C2").Select()
